# ALZ Policy Assignments workbook update
#
# 1) Rename the first worksheet from "ALZ Policy Assignments 03CY23" to
#    "ALZ Policy Assignments 12CY23". The workbook-scoped hidden defined
#    name "_xlnm._FilterDatabase" refers to this sheet by name, so Excel
#    updates that reference automatically when the sheet is renamed.
# 2) Reset the sheet's view: scroll back to the top-left (dropping the old
#    topLeftCell="A6" scroll position) and move the active selection from
#    K16 to A2.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALZ Policy Assignments 03CY23")
$ws.Name = "ALZ Policy Assignments 12CY23"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2").Select()
